$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 1 & 2 price updates (~15% increase)
$ws.Range("D22").Value = 378.388
$ws.Range("D23").Value = 479.507
$ws.Range("D34").Value = 373.17
$ws.Range("D35").Value = 515.715
$ws.Range("D45").Value = 457.978
$ws.Range("D46").Value = 513.7569999999999
